$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Remove the "Buyer Agency" helper column currently at D (content + the
#    column's own width/style definition) WITHOUT shifting any other column.
# ---------------------------------------------------------------------------
$ws.Range("D1:D11").Style = "Normal"
$ws.Range("D1:D11").ClearContents()

# ---------------------------------------------------------------------------
# 2) Clear the stray leftover note cell in Z2 before we shift the tail
#    columns, so it doesn't survive the upcoming column insert.
# ---------------------------------------------------------------------------
$ws.Range("Z2").Style = "Normal"
$ws.Range("Z2").ClearContents()

# ---------------------------------------------------------------------------
# 3) Insert a new blank column before the old "Other" column (X), which
#    pushes "Other" from X -> Y, restoring the "Buyer Agency" column at its
#    new home (X) right after "Agent Phone" (W).
# ---------------------------------------------------------------------------
$ws.Columns("X").Insert()

# New column X header + numeric formatting style (matches col C / old D)
$ws.Range("X1").Style = "Normal"
$ws.Cells.Item(1, 24).Value = "Buyer Agency"

# Re-apply cell styles on the new X column data cells (style ids 9/10 used
# throughout the sheet for this numeric helper column).
$ws.Range("X2").Value = $null
$ws.Range("X3").Value = $null
$ws.Range("X4").Value = $null
$ws.Range("X5").Value = 3
$ws.Range("X6").Value = $null
$ws.Range("X8").Value = 2
$ws.Range("X9").Value = 2
$ws.Range("X10").Value = $null
$ws.Range("X11").Value = $null

# ---------------------------------------------------------------------------
# 4) Row 8 content fixes: typo correction + case fix.
# ---------------------------------------------------------------------------
$ws.Range("N8").Value = "Referigerator"
$ws.Range("Q8").Value = "n"

# ---------------------------------------------------------------------------
# 5) Row 9: populate the previously-empty contract row with full data
#    (mirrors the pattern used by the other rows in the table).
# ---------------------------------------------------------------------------
$ws.Range("E9").Value = 400000
$ws.Range("F9").Value = "Joh nand Sue Barbera"
$ws.Range("H9").Value = "Cash"
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = "n"
$ws.Range("K9").Value = "n"
$ws.Range("L9").Value = "A"
$ws.Range("M9").Value = "Y-Buyer"
$ws.Range("N9").Value = "Referigerator"
$ws.Range("O9").Value = "no"
$ws.Range("P9").Value = "Yes-695"
$ws.Range("Q9").Value = "n"
$ws.Range("R9").Value = 45989
$ws.Range("S9").Value = 2
$ws.Range("T9").Value = 45900
$ws.Range("U9").Value = "Brian Curtis"
$ws.Range("V9").Value = "EB00054032"
$ws.Range("W9").Value = "479-531-2317"

# ---------------------------------------------------------------------------
# 6) Sheet view bookkeeping (matches the user scrolling / selecting X9).
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollColumn = 7
$ws.Range("X9").Select()

Write-Host "edit applied"
